# Update the "Förändrad" (Changed) date column (C) for every data row
# that currently holds the serial date 46075 (2026-02-22), bumping it
# forward one day to 46076 (2026-02-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 239 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value = 46076
    }
}
